$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 12:03"

# --- Row 18 (rank 22, Banglades): refreshed case counts ---
$ws.Range("B18").Value = 304583
$ws.Range("C18").Value = 2436
$ws.Range("D18").Value = 193458
$ws.Range("E18").Value = 106998
$ws.Range("G18").Value = 45
$ws.Range("H18").Value = 4127

# --- Rows 47/48: Polonia and Japon swap rank positions with refreshed data ---
# Row 47 was Japon, now becomes Polonia (with new data)
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 64689
$ws.Range("C47").Value = 887
$ws.Range("D47").Value = 44097
$ws.Range("E47").Value = 18582
$ws.Range("G47").Value = 16
$ws.Range("H47").Value = 2010

# Row 48 was Polonia, now becomes Japon (carrying Japon's former data)
$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 63822
$ws.Range("D48").Value = 51688
$ws.Range("E48").Value = 10925
$ws.Range("H48").Value = 1209

# --- Row 92 (rank 96, Malasia): refreshed case counts ---
$ws.Range("B92").Value = 9296
$ws.Range("C92").Value = 5
$ws.Range("D92").Value = 8994
$ws.Range("E92").Value = 177

# --- Row 101 (rank 105, Finlandia): refreshed case counts ---
$ws.Range("B101").Value = 8019
$ws.Range("C101").Value = 17
$ws.Range("E101").Value = 484

# --- Row 111 (rank 115, Hong Kong): refreshed case counts ---
$ws.Range("B111").Value = 4756
$ws.Range("C111").Value = 20
$ws.Range("D111").Value = 4200
$ws.Range("E111").Value = 475
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = 81

# --- Row 194 (rank 198, Liechtenstein): refreshed case counts ---
$ws.Range("B194").Value = 105
$ws.Range("C194").Value = 3
$ws.Range("D194").Value = 96
$ws.Range("E194").Value = 8
